# Update "想去人数" (number of people wanting to attend) figures that changed
# between scrapes, on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 1470
    $ws.Range("F3").Value = 3090
    $ws.Range("F5").Value = 691
}
